# Scheduled-runner style refresh of market/profit figures across the
# Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW). Each block below
# re-points currentAveragePrice / LevePrice / LeveProfit columns (H:N) for
# the rows whose sourced prices moved since the last run. A few rows also
# gain/lose a profit cell (M/N) because HQ or NQ pricing data became
# available/unavailable for that leve on this pass.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 25.25
$ws.Range("I9").Value = 40
$ws.Range("K9").Value = 40
$ws.Range("M9").Value = 129
$ws.Range("H15").Value = 2102.375
$ws.Range("I15").Value = 2102.375
$ws.Range("K15").Value = 6307.125
$ws.Range("M15").Value = -6138.125
$ws.Range("H17").Value = 11896
$ws.Range("J17").Value = 12995.556
$ws.Range("L17").Value = 38986.66800000001
$ws.Range("N17").Value = -39322.66800000001
$ws.Range("H51").Value = 9814.125
$ws.Range("I51").Value = 12952.75
$ws.Range("J51").Value = 6675.5
$ws.Range("K51").Value = 12952.75
$ws.Range("L51").Value = 6675.5
$ws.Range("M51").Value = -12468.75
$ws.Range("N51").Value = -7643.5
$ws.Range("H76").Value = 4996.3335
$ws.Range("I76").Value = 4996.3335
$ws.Range("K76").Value = 4996.3335
$ws.Range("M76").Value = -4681.3335
$ws.Range("H79").Value = 4996.3335
$ws.Range("I79").Value = 4996.3335
$ws.Range("K79").Value = 4996.3335
$ws.Range("M79").Value = -3904.3335
$ws.Range("H86").Value = 169682.67
$ws.Range("I86").Value = 253424.25
$ws.Range("K86").Value = 253424.25
$ws.Range("M86").Value = -252301.25
$ws.Range("H89").Value = 169682.67
$ws.Range("I89").Value = 253424.25
$ws.Range("K89").Value = 1267121.25
$ws.Range("M89").Value = -1261505.25
$ws.Range("H116").Value = 4866.212
$ws.Range("I116").Value = 4872.05
$ws.Range("J116").Value = 4857.231
$ws.Range("K116").Value = 4872.05
$ws.Range("L116").Value = 4857.231
$ws.Range("M116").Value = -1430.05
$ws.Range("N116").Value = -11741.231
$ws.Range("H118").Value = 920.9583
$ws.Range("I118").Value = 479.625
$ws.Range("K118").Value = 1438.875
$ws.Range("M118").Value = 218.125
$ws.Range("H137").Value = 2353.9412
$ws.Range("I137").Value = 1051.4286
$ws.Range("J137").Value = 3265.7
$ws.Range("K137").Value = 3154.2858
$ws.Range("L137").Value = 9797.099999999999
$ws.Range("M137").Value = -604.2857999999997
$ws.Range("N137").Value = -14897.1
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36183.305
$ws.Range("I32").Value = 39157.035
$ws.Range("K32").Value = 39157.035
$ws.Range("M32").Value = -38870.035
$ws.Range("H45").Value = 5590
$ws.Range("I45").Value = 3512
$ws.Range("J45").Value = 5886.857
$ws.Range("K45").Value = 3512
$ws.Range("L45").Value = 5886.857
$ws.Range("M45").Value = -3135
$ws.Range("N45").Value = -6640.857
$ws.Range("H61").Value = 2191.8
$ws.Range("I61").Value = 2157.5557
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 2157.5557
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -1945.5557
$ws.Range("N61").Value = -2924
$ws.Range("H74").Value = 63415.75
$ws.Range("I74").Value = 63415.75
$ws.Range("K74").Value = 63415.75
$ws.Range("M74").Value = -62541.75
$ws.Range("H77").Value = 63415.75
$ws.Range("I77").Value = 63415.75
$ws.Range("K77").Value = 317078.75
$ws.Range("M77").Value = -312710.75
$ws.Range("H82").Value = 40000
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40722
$ws.Range("H85").Value = 40000
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42496
$ws.Range("H94").Value = 29285.715
$ws.Range("J94").Value = 29285.715
$ws.Range("L94").Value = 29285.715
$ws.Range("N94").Value = -31087.715
$ws.Range("H102").Value = 3726.923
$ws.Range("I102").Value = 3329.6
$ws.Range("K102").Value = 3329.6
$ws.Range("M102").Value = -1707.6
$ws.Range("H131").Value = 120000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 120000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 120000
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -130080
$ws.Range("H136").Value = 2191.8
$ws.Range("I136").Value = 2157.5557
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 6472.6671
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -3922.6671
$ws.Range("N136").Value = -12600
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 5179.2
$ws.Range("I64").Value = 2499.5
$ws.Range("J64").Value = 6965.6665
$ws.Range("K64").Value = 2499.5
$ws.Range("L64").Value = 6965.6665
$ws.Range("M64").Value = -2274.5
$ws.Range("N64").Value = -7415.6665
$ws.Range("H67").Value = 5179.2
$ws.Range("I67").Value = 2499.5
$ws.Range("J67").Value = 6965.6665
$ws.Range("K67").Value = 2499.5
$ws.Range("L67").Value = 6965.6665
$ws.Range("M67").Value = -1719.5
$ws.Range("N67").Value = -8525.666499999999
$ws.Range("H86").Value = 2270.0908
$ws.Range("I86").Value = 2121.375
$ws.Range("J86").Value = 2666.6667
$ws.Range("K86").Value = 2121.375
$ws.Range("L86").Value = 2666.6667
$ws.Range("M86").Value = -998.375
$ws.Range("N86").Value = -4912.6667
$ws.Range("H89").Value = 2270.0908
$ws.Range("I89").Value = 2121.375
$ws.Range("J89").Value = 2666.6667
$ws.Range("K89").Value = 10606.875
$ws.Range("L89").Value = 13333.3335
$ws.Range("M89").Value = -4990.875
$ws.Range("N89").Value = -24565.3335
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2350.68
$ws.Range("I31").Value = 1792.8889
$ws.Range("K31").Value = 1792.8889
$ws.Range("M31").Value = -1497.8889
$ws.Range("H34").Value = 2350.68
$ws.Range("I34").Value = 1792.8889
$ws.Range("K34").Value = 1792.8889
$ws.Range("M34").Value = -1590.8889
$ws.Range("H134").Value = 114618.78
$ws.Range("I134").Value = 168259.17
$ws.Range("K134").Value = 504777.51
$ws.Range("M134").Value = -502242.51
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1624
$ws.Range("I5").Value = 1582.6666
$ws.Range("J5").Value = 1748
$ws.Range("K5").Value = 4747.9998
$ws.Range("L5").Value = 5244
$ws.Range("M5").Value = -4635.9998
$ws.Range("N5").Value = -5468
$ws.Range("H61").Value = 307.5
$ws.Range("J61").Value = 425
$ws.Range("L61").Value = 1275
$ws.Range("N61").Value = -1705
$ws.Range("H135").Value = 1624
$ws.Range("I135").Value = 1582.6666
$ws.Range("J135").Value = 1748
$ws.Range("K135").Value = 14243.9994
$ws.Range("L135").Value = 15732
$ws.Range("M135").Value = -11708.9994
$ws.Range("N135").Value = -20802
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 54034.26
$ws.Range("I107").Value = 84417.414
$ws.Range("K107").Value = 84417.414
$ws.Range("M107").Value = -82497.414
$ws.Range("H113").Value = 2772.3572
$ws.Range("I113").Value = 2198.5
$ws.Range("J113").Value = 4207
$ws.Range("K113").Value = 2198.5
$ws.Range("L113").Value = 4207
$ws.Range("M113").Value = -28.5
$ws.Range("N113").Value = -8547
$ws.Range("H132").Value = 28927.541
$ws.Range("I132").Value = 36089.69
$ws.Range("J132").Value = 2964.75
$ws.Range("K132").Value = 108269.07
$ws.Range("L132").Value = 8894.25
$ws.Range("M132").Value = -105739.07
$ws.Range("N132").Value = -13954.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 36574.75
$ws.Range("I22").Value = 62219.11
$ws.Range("J22").Value = 3603.4285
$ws.Range("K22").Value = 62219.11
$ws.Range("L22").Value = 3603.4285
$ws.Range("M22").Value = -61924.11
$ws.Range("N22").Value = -4193.4285
$ws.Range("H27").Value = 36574.75
$ws.Range("I27").Value = 62219.11
$ws.Range("J27").Value = 3603.4285
$ws.Range("K27").Value = 62219.11
$ws.Range("L27").Value = 3603.4285
$ws.Range("M27").Value = -62112.11
$ws.Range("N27").Value = -3817.4285
$ws.Range("H46").Value = 12766
$ws.Range("I46").Value = 17383.55
$ws.Range("J46").Value = 5070.0835
$ws.Range("K46").Value = 17383.55
$ws.Range("L46").Value = 5070.0835
$ws.Range("M46").Value = -17195.55
$ws.Range("N46").Value = -5446.0835
$ws.Range("H82").Value = 2113.2334
$ws.Range("J82").Value = 2465.5
$ws.Range("L82").Value = 2465.5
$ws.Range("N82").Value = -3187.5
$ws.Range("H85").Value = 2113.2334
$ws.Range("J85").Value = 2465.5
$ws.Range("L85").Value = 2465.5
$ws.Range("N85").Value = -4961.5
$ws.Range("H136").Value = 8358.5
$ws.Range("I136").Value = 7811.3335
$ws.Range("K136").Value = 23434.0005
$ws.Range("M136").Value = -20884.0005
